# "Add strength to attacks"
# Rebalance the per-level stat tables: raise the "strength" (column D) values
# and lower the matching "attack" stat (intelligence/wisdom/intuition) that
# used to carry all the extra power on levels 1 and 2, so some of that power
# now comes from strength instead.

$wb = $excel.ActiveWorkbook

# --- tough_levels: view/selection only, no stat changes -------------------
$wsTough = $wb.Worksheets.Item("tough_levels")
$wsTough.Range("D8").Select()

# --- intelligent_levels: strength up, intelligence down -------------------
$wsIntelligent = $wb.Worksheets.Item("intelligent_levels")
$wsIntelligent.Range("D2").Value = 10
$wsIntelligent.Range("G2").Value = 10
$wsIntelligent.Range("D3").Value = 15
$wsIntelligent.Range("G3").Value = 25
$wsIntelligent.Range("G8").Select()

# --- wise_levels: strength up, wisdom down ---------------------------------
$wsWise = $wb.Worksheets.Item("wise_levels")
$wsWise.Range("D2").Value = 10
$wsWise.Range("E2").Value = 10
$wsWise.Range("D3").Value = 15
$wsWise.Range("E3").Value = 25
$wsWise.Range("E11").Select()

# --- intuitive_levels: strength up, intuition down; this ends as the ------
# --- active sheet/tab, matching the final saved workbook state ------------
$wsIntuitive = $wb.Worksheets.Item("intuitive_levels")
$wsIntuitive.Range("D2").Value = 10
$wsIntuitive.Range("F2").Value = 10
$wsIntuitive.Range("D3").Value = 15
$wsIntuitive.Range("F3").Value = 25
$wsIntuitive.Activate()
$wsIntuitive.Range("F12").Select()
